$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "https://careers.homedepot.com/"
$ws.Range("C2").Value = "https://corporate.homedepot.com/"
$ws.Range("C3").Value = "https://corporate.homedepot.com/newsroom"
$ws.Range("C4").Value = "https://corporate.homedepot.com/foundation"
$ws.Range("C5").Value = "https://ir.homedepot.com/"
$ws.Range("C6").Value = "/c/Government_Customers"
$ws.Range("C7").Value = "/c/suppliers_and_providers"
$ws.Range("C8").Value = "https://www.homedepot.com/c/SF_MS_The_Home_Depot_Affiliate_Program"
$ws.Range("C9").Value = "https://ecoactions.homedepot.com/"
